# Simulation for sphinx3 added
# Fills in the previously-empty sphinx3 rows (71-74) on both the Config1
# and Config2 sheets with simulation results, and updates the
# window/selection state so that Config2 ends up the active/selected sheet.

$wb = $excel.ActiveWorkbook

$config1 = $wb.Worksheets.Item("Config1")
$config2 = $wb.Worksheets.Item("Config2")

# ---------------------------------------------------------------------
# Config1 sheet - sphinx3 results (rows 71-74: LRU, SRRIP, Hawkeye, OPTGen)
# ---------------------------------------------------------------------

# Row 71 - LRU
$config1.Range("C71").Value = 50000000
$config1.Range("D71").Value = 92768345
$config1.Range("E71").Value = 754197
$config1.Range("F71").Value = 160196
$config1.Range("G71").Value = 594001
$config1.Range("H71").Formula = "=(C71/D71)"
$config1.Range("I71").Formula = "=G71/(C71/1000)"

# Row 72 - SRRIP
$config1.Range("C72").Value = 50000000
$config1.Range("D72").Value = 89747896
$config1.Range("E72").Value = 754200
$config1.Range("F72").Value = 215084
$config1.Range("G72").Value = 539116
$config1.Range("H72").Formula = "=(C72/D72)"
$config1.Range("I72").Formula = "=G72/(C72/1000)"

# Row 73 - Hawkeye
$config1.Range("C73").Value = 50000000
$config1.Range("D73").Value = 78116641
$config1.Range("E73").Value = 754196
$config1.Range("F73").Value = 376654
$config1.Range("G73").Value = 377542
$config1.Range("H73").Formula = "=(C73/D73)"
$config1.Range("I73").Formula = "=G73/(C73/1000)"

# Row 74 - OPTGen
$config1.Range("C74").Value = 50000000
$config1.Range("D74").Value = 78116641
$config1.Range("E74").Value = 22436
$config1.Range("F74").Value = 11343
$config1.Range("G74").Formula = "=E74-F74"
$config1.Range("H74").Formula = "=(C74/D74)"
$config1.Range("I74").Formula = "=G74/(C74/1000)"
$config1.Range("J74").Formula = "=F74/E74"

# ---------------------------------------------------------------------
# Config2 sheet - sphinx3 results (rows 71-74: LRU, SRRIP, Hawkeye, OPTGen)
# ---------------------------------------------------------------------

# Row 71 - LRU
$config2.Range("C71").Value = 50000000
$config2.Range("D71").Value = 60373750
$config2.Range("E71").Value = 836292
$config2.Range("F71").Value = 153218
$config2.Range("G71").Value = 683074
$config2.Range("H71").Formula = "=(C71/D71)"
$config2.Range("I71").Formula = "=G71/(C71/1000)"

# Row 72 - SRRIP
$config2.Range("C72").Value = 50000000
$config2.Range("D72").Value = 58442276
$config2.Range("E72").Value = 836345
$config2.Range("F72").Value = 215632
$config2.Range("G72").Value = 620713
$config2.Range("H72").Formula = "=(C72/D72)"
$config2.Range("I72").Formula = "=G72/(C72/1000)"

# Row 73 - Hawkeye
$config2.Range("C73").Value = 50000000
$config2.Range("D73").Value = 53262804
$config2.Range("E73").Value = 836625
$config2.Range("F73").Value = 309808
$config2.Range("G73").Value = 526817
$config2.Range("H73").Formula = "=(C73/D73)"
$config2.Range("I73").Formula = "=G73/(C73/1000)"

# Row 74 - OPTGen
$config2.Range("C74").Value = 50000000
$config2.Range("D74").Value = 53262804
$config2.Range("E74").Value = 6303
$config2.Range("F74").Value = 5767
$config2.Range("G74").Formula = "=E74-F74"
$config2.Range("H74").Formula = "=(C74/D74)"
$config2.Range("I74").Formula = "=G74/(C74/1000)"
$config2.Range("J74").Formula = "=F74/E74"

# ---------------------------------------------------------------------
# View state - scroll/selection on each sheet, and make Config2 active
# ---------------------------------------------------------------------

$config1.Application.Goto($config1.Range("A55"))
$config1.Range("D77").Select()

$config2.Application.Goto($config2.Range("A61"))
$config2.Range("E75").Select()

$config2.Activate()
